$d = $word.ActiveDocument

# --- Step 1: remove the four paragraphs that were dropped entirely ---
# (the CPAB-derivation / GeLU-variant / regularization / closing-remark paragraphs)
$prefixesToRemove = @(
    "המאמר מתבונן במקרה של שדה וקטורי",
    "המאמר מציע לשכלל את פונקציית אקטיבציה",
    "בנוסף יש איבר רגולריזציה",
    "מאמר כיפי וכתוב היטב"
)

$indicesToRemove = @()
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($pref in $prefixesToRemove) {
        if ($t.StartsWith($pref)) {
            $indicesToRemove += $i
        }
    }
    $i += 1
}
$sortedDesc = $indicesToRemove | Sort-Object -Descending
foreach ($idx in $sortedDesc) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# --- Step 2: text replacements ---
# date in the heading (18.07.24 -> 17.07.24)
$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק 18.07.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק 17.07.24: ⚡️🚀", 2) | Out-Null

# title (Trainable Highly-expressive Activation Functions -> Learning Rate Curriculum)
$d.Content.Find.Execute("Trainable Highly-expressive Activation Functions", $true, $false, $false, $false, $false, $true, 1, $false, "Learning Rate Curriculum", 2) | Out-Null

# intro paragraph
$d.Content.Find.Execute(" ממשיכים את קו הגיוון וסוקרים מאמר לא קשור ישירות למודלי שפה. היום נסקור מאמר של כמה חוקרים ישראלים המציע דרך חדשה לבנות פונקציות אקטיבציה ברשת נוירונים. היום פונקציות אקטיבציה הן לא נלמדות לרוב (ReLU, GeLU, tanh וכדומה). לפעמים פונקציות אקטיבציה מכילות hyperparameter שלא נלמד במהלך האימון אלא נקבע מראש (Leaky ReLU, Swish וכדומה).", $true, $false, $false, $false, $false, $true, 1, $false, "רוב המאמרים שסקרתי לאחרונה היו בנושא מודלי שפה והחלטתי לגוון טיפה ולסקור מאמרים בנושאים אחרים. מאמר שנסקור היום מדבר על שיטת אימון הנקראת למידת curriculum שבא אנו מאמנים את המודל כמו שאנו מלמדים חומר לתלמידים - מהקל לקשה. יש כמה וריאציות של למידת curriculum: באחת מהם אנו מתחילים לאמן מודל עם דוגמאות קלות ובהדרגה מעלים את קושי הדוגמאות. הוריאציה השניה אנו מתחילים ממשימה קלה יותר ומעלים את מורכבותה בהדרגה. בשלישית מאמנים מודל יחסית פשוט ומעלים את ״מורכבות״ של המודל. ", 2) | Out-Null

# paragraph describing the curriculum idea for learning-rate
$d.Content.Find.Execute("המאמר מציע פונקציות אקטיבציה שהן(הפרמטרים שלהן) אשכרה נלמדות במהלך האימון. ד״א לאחרונה ראינו דוגמא נוספת לפונקצית אקטיבציה נלמדת ראינו לא מזמן במאמר המפורסם Kolmogorov-Arnold network או KAN - שם אלו היו ספליינים נלמדים. במאמר המסוקר אימצו שיטה אחרת לבנייה של פונקציות אקטיבציה נלמדות. הבנייה נעשה דרך שדות וקטורים שמגדירות את המסלול של נקודה במרחב. ", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר מציע גישת curriculum אבל לקצב למידה. המחברים מציינים שלמשל ברשתות קונבולוציה עדיף בהתחלה להתמקד יותר בלמידה של השכבות הראשונות כי למעשה אם אלו לא נלמדו טוב ועדיין קרובים למצב האיתחול שלהם אז הם יוצרים דאטה ״רועש, מדי שזורם גם לשכבות הבאות שמתקשות להתמודד איתו (המאמר מציין כמה עבודות שחקרו את הנושא והגיעו למסקנות האלו). תופעה דומה מתרחשת גם כאשר אנו עושים פיין טיון למודל למשימה מסוימת כאשר המודל לפני זה אומן למשימה אחרת. ", 2) | Out-Null

# paragraph describing how authors address the issue
$d.Content.Find.Execute("במקרה הזה אנו מתחילים מנקודה x ובעזרת נגזרת של כיוון תנועת הנקודה(=שדה וקטורי) ב״זמן״ (שמתחיל ב t=0 ומסתיים ב- t=1)  נבנה המסלול של נקודה x. המסלול מסתיים ב t=1 לכל x שלמעשה מגדיר לנו פונקציית אקטיבציה (a(x. ניתן לתאר את התקדמות נקודה באמצעות משוואה אינטגרלית (כמו שיטת אוילר לפתרון משוואות דיפרנציאליות).", $true, $false, $false, $false, $false, $true, 1, $false, "כדי להתמודד עם סוגיה זו המחברים מציעים להתחיל מקצב למידה גבוה עבור השכבות הראשונות (שיורד ככל שמתקדמים לשכבות עמוקות יותר). במהלך האיטרציות לעלות את קצב למידה בכל השכבות כך (קצב עלייה לא שווה בין השכבות) כך שעם הזמן (=איטרציות) קצבי הלמידה של כל השכבות משתוות. נציין שהמחברים מציעים שמספר האיטרציות הנדרש להשוואת קצב הלמידה עבור כל השכבות צריך להיות משמעותית קטן יותר מכמות האיטרציות הכולל הנדרש לאימון המודל. כלומר כל השיטה הזו מופעלת בשלב ה״חימום״ של הרשת. ", 2) | Out-Null

# arxiv link
$d.Content.Find.Execute("https://arxiv.org/abs/2407.07564", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2205.09180", 2) | Out-Null

